# Applies the "Nexus 5X - Magnetic Field - Raw" update:
#  - replaces the magnetic-field sample readings in A2:A31
#  - updates the active selection / scroll position of the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    117.82044,
    124.62372000000001,
    119.985119999999,
    119.8305,
    120.44898000000001,
    120.6036,
    124.77834,
    120.44898000000001,
    120.13974,
    120.758219999999,
    119.8305,
    121.22208000000001,
    119.36664,
    120.13974,
    119.36664,
    119.36664,
    123.077519999999,
    118.12967999999999,
    120.758219999999,
    122.1498,
    118.28429999999901,
    120.758219999999,
    120.91284,
    119.67588000000001,
    120.44898000000001,
    121.995179999999,
    123.54138,
    119.36664,
    115.34652,
    119.8305
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Move the view so row 22 is at the top and D29 is the active/selected cell,
# matching the sheetView/selection state recorded in the workbook.
$win = $excel.ActiveWindow
$ws.Range("D29").Select()
$win.ScrollRow = 22
$win.ScrollColumn = 1
